$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.693.62'
$ws.Range('E2').Value = '  -3.65%  '
$ws.Range('D3').Value = '2.561.04'
$ws.Range('E3').Value = '  -1.58%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '''507.92'
$ws.Range('E5').Value = '  -2.84%  '
$ws.Range('D6').Value = '''143.60'
$ws.Range('E6').Value = '  -7.24%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '''0.556'
$ws.Range('E8').Value = '  -6.15%  '
$ws.Range('D9').Value = '2.562.53'
$ws.Range('E9').Value = '  -1.73%  '
$ws.Range('D10').Value = '''6.21'
$ws.Range('E10').Value = '  -7.28%  '
$ws.Range('D11').Value = '''0.102'
$ws.Range('E11').Value = '  -3.53%  '
$ws.Range('D12').Value = '''0.331'
$ws.Range('E12').Value = '  -4.66%  '
$ws.Range('E13').Value = '  -1.00%  '
$ws.Range('D14').Value = '3.004.52'
$ws.Range('E14').Value = '  -1.81%  '
$ws.Range('D15').Value = '58.650.99'
$ws.Range('E15').Value = '  -3.78%  '
$ws.Range('D16').Value = '''20.59'
$ws.Range('E16').Value = '  -5.21%  '
$ws.Range('E17').Value = '  -4.91%  '
$ws.Range('D18').Value = '2.555.93'
$ws.Range('E18').Value = '  -1.78%  '
$ws.Range('D19').Value = '''4.52'
$ws.Range('E19').Value = '  -5.24%  '
$ws.Range('D20').Value = '''333.39'
$ws.Range('E20').Value = '  -5.86%  '
$ws.Range('D21').Value = '''10.06'
$ws.Range('E21').Value = '  -4.82%  '
$ws.Range('D22').Value = '''0.997'
$ws.Range('E22').Value = '  -0.21%  '
$ws.Range('D23').Value = '''5.95'
$ws.Range('E23').Value = '  -4.37%  '
$ws.Range('D24').Value = '''59.89'
$ws.Range('E24').Value = '  -1.58%  '
$ws.Range('D25').Value = '''0.407'
$ws.Range('E25').Value = '  -4.65%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('E27').Value = '  -5.93%  '
$ws.Range('D28').Value = '0.0₃0781'
$ws.Range('E28').Value = '  -7.96%  '
$ws.Range('D29').Value = '''6.90'
$ws.Range('E29').Value = '  -7.03%  '
$ws.Range('E30').Value = '  -0.04%  '
$ws.Range('D31').Value = '''149.51'
$ws.Range('E31').Value = '  +0.99%  '
$ws.Range('D32').Value = '''5.85'
$ws.Range('E32').Value = '  -7.01%  '
$ws.Range('D33').Value = '''18.55'
$ws.Range('E33').Value = '  -4.36%  '
$ws.Range('E34').Value = '  -3.66%  '
$ws.Range('D35').Value = '''3.95'
$ws.Range('E35').Value = '  -5.78%  '
$ws.Range('D36').Value = '''0.900'
$ws.Range('E36').Value = '  -2.18%  '
$ws.Range('E37').Value = '  -8.13%  '
$ws.Range('D38').Value = '''35.94'
$ws.Range('E38').Value = '  -1.49%  '
$ws.Range('D39').Value = '''0.822'
$ws.Range('E39').Value = '  -6.26%  '
$ws.Range('D40').Value = '''288.54'
$ws.Range('E40').Value = '  -0.57%  '
$ws.Range('E41').Value = '  -8.34%  '
$ws.Range('D42').Value = '''3.53'
$ws.Range('E42').Value = '  -7.36%  '
$ws.Range('D43').Value = '''0.997'
$ws.Range('E43').Value = '  -0.01%  '
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D44').Value = '''0.0981'
$ws.Range('E44').Value = '  -3.53%  '
$ws.Range('B45').Value = 'Mantle'
$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D45').Value = '''0.606'
$ws.Range('E45').Value = '  -2.56%  '
$ws.Range('D46').Value = '''0.0533'
$ws.Range('E46').Value = '  -4.98%  '
$ws.Range('D47').Value = '''18.70'
$ws.Range('E47').Value = '  -4.47%  '
$ws.Range('E48').Value = '  -0.09%  '
$ws.Range('E49').Value = '  -4.59%  '
$ws.Range('E50').Value = '  -7.97%  '
$ws.Range('D51').Value = '1.912.61'
$ws.Range('E51').Value = '  -2.71%  '
